# Apply the changes described by the commit:
# "Added new charts.  Added function and calculation of IRRs for the trials."
#
# Concretely (visible in the data itself) this adds one new simulation
# scenario ("Shiny1") to Sheet1 and one matching cash-flow / withdrawal
# row ("Withdrawals") to the cashflows sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Sheet1
$ws2 = $wb.Worksheets.Item(2)   # cashflows

# ---------------------------------------------------------------
# Sheet1: new row 8 - "Shiny1" simulation definition
# ---------------------------------------------------------------
$ws1.Range("A8").Value = "Shiny1"
$ws1.Range("B8").Value = 500
$ws1.Range("C8").Value = 1000000
$ws1.Range("D8").Value = "F"
$ws1.Range("E8").Value = 30
$ws1.Range("F8").Value = 101
$ws1.Range("G8").Value = 0.025
$ws1.Range("H8").Value = 0.08
$ws1.Range("I8").Value = 0.09
$ws1.Range("J8").Value = 0.01
$ws1.Range("K8").Value = $false
$ws1.Range("L8").Value = 0.6
$ws1.Range("M8").Value = 12
$ws1.Range("N8").Value = 0
$ws1.Range("Q8").Value = $false

# R8/S8 reuse the same date / returnGeneratorMethod style (s="1") that the
# rest of the column already uses, so copy formatting from the row above
# before writing the new values.
$ws1.Range("R7").Copy($ws1.Range("R8"))
$ws1.Range("R8").Formula = "=TODAY()"

$ws1.Range("S7").Copy($ws1.Range("S8"))
$ws1.Range("S8").Value = "S"

$ws1.Range("T8").Select()

# ---------------------------------------------------------------
# cashflows: new row 20 - "Withdrawals" cash flow tied to "Shiny1"
# ---------------------------------------------------------------
$ws2.Range("A20").Value = "Shiny1"
$ws2.Range("B20").Value = "Withdrawals"
$ws2.Range("C20").Value = "start"
$ws2.Range("D20").Value = 1
$ws2.Range("E20").Value = "end"
$ws2.Range("F20").Value = 30
$ws2.Range("G20").Value = "w"
$ws2.Range("H20").Value = 40000
$ws2.Range("I20").Value = $true
$ws2.Range("J20").Value = 0

$ws2.Activate()
$ws2.Range("J20").Select()
